# Finishing touches to tutorial: fix values in the Consumption sheet and
# update the active sheet/selection, mirroring the author's Excel session.

$wb = $excel.ActiveWorkbook

# --- Consumption sheet: correct a handful of consumption-matrix values ---
$consumption = $wb.Worksheets.Item("Consumption")
$consumption.Range("C3").Value = 0.99
$consumption.Range("B4").Value = 0.48
$consumption.Range("C4").Value = 0.49
$consumption.Range("B5").Value = 0.49
$consumption.Range("C5").Value = 0.48

# Make Consumption the active sheet (was Mix), with C5 as the active/selected cell
$consumption.Activate()
$consumption.Range("C5").Select() | Out-Null
